# Auto-generated edit script applying the Sagittarius_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 65000
$ws.Range("J3").Value = 65000
$ws.Range("L3").Value = 65000
$ws.Range("N3").Value = -65228
$ws.Range("H17").Value = 7494.2856
$ws.Range("J17").Value = 7494.2856
$ws.Range("L17").Value = 22482.8568
$ws.Range("N17").Value = -22818.8568
$ws.Range("H74").Value = 75370.336
$ws.Range("I74").Value = 142974.6
$ws.Range("K74").Value = 142974.6
$ws.Range("M74").Value = -142038.6
$ws.Range("H77").Value = 75370.336
$ws.Range("I77").Value = 142974.6
$ws.Range("K77").Value = 714873
$ws.Range("M77").Value = -710193
$ws.Range("H98").Value = 3451.5264
$ws.Range("I98").Value = 975.9
$ws.Range("J98").Value = 6202.222
$ws.Range("K98").Value = 975.9
$ws.Range("L98").Value = 6202.222
$ws.Range("M98").Value = 522.1
$ws.Range("N98").Value = -9198.222
$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490
$ws.Range("H122").Value = 3451.5264
$ws.Range("I122").Value = 975.9
$ws.Range("J122").Value = 6202.222
$ws.Range("K122").Value = 2927.7
$ws.Range("L122").Value = 18606.666
$ws.Range("M122").Value = -477.6999999999998
$ws.Range("N122").Value = -23506.666
$ws.Range("H132").Value = 1505
$ws.Range("I132").Value = 1507.8334
$ws.Range("K132").Value = 4523.5002
$ws.Range("M132").Value = -1993.5002
$ws.Range("H135").Value = 620.8946999999999
$ws.Range("I135").Value = 520.3889
$ws.Range("J135").Value = 2430
$ws.Range("K135").Value = 4683.5001
$ws.Range("L135").Value = 21870
$ws.Range("M135").Value = -2148.5001
$ws.Range("N135").Value = -26940
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2063.5264
$ws.Range("I74").Value = 1569.75
$ws.Range("K74").Value = 1569.75
$ws.Range("M74").Value = -695.75
$ws.Range("H77").Value = 2063.5264
$ws.Range("I77").Value = 1569.75
$ws.Range("K77").Value = 7848.75
$ws.Range("M77").Value = -3480.75
$ws.Range("H134").Value = 81663.336
$ws.Range("J134").Value = 81663.336
$ws.Range("L134").Value = 81663.336
$ws.Range("N134").Value = -91803.336
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 682
$ws.Range("I20").Value = 682
$ws.Range("K20").Value = 682
$ws.Range("M20").Value = -435
$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680
$ws.Range("H134").Value = 4587
$ws.Range("I134").Value = 4587
$ws.Range("K134").Value = 13761
$ws.Range("M134").Value = -11226
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1691.8
$ws.Range("I16").Value = 1129.5
$ws.Range("J16").Value = 3003.8333
$ws.Range("K16").Value = 1129.5
$ws.Range("L16").Value = 3003.8333
$ws.Range("M16").Value = -842.5
$ws.Range("N16").Value = -3577.8333
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""
$ws.Range("H107").Value = 755.0714
$ws.Range("I107").Value = 414
$ws.Range("J107").Value = 1369
$ws.Range("K107").Value = 414
$ws.Range("L107").Value = 1369
$ws.Range("M107").Value = 1506
$ws.Range("N107").Value = -5209
$ws.Range("H113").Value = 1691.8
$ws.Range("I113").Value = 1129.5
$ws.Range("J113").Value = 3003.8333
$ws.Range("K113").Value = 1129.5
$ws.Range("L113").Value = 3003.8333
$ws.Range("M113").Value = 1040.5
$ws.Range("N113").Value = -7343.8333
$ws.Range("H134").Value = 2671.2
$ws.Range("I134").Value = 2637.5386
$ws.Range("J134").Value = 2890
$ws.Range("K134").Value = 7912.6158
$ws.Range("L134").Value = 8670
$ws.Range("M134").Value = -5377.6158
$ws.Range("N134").Value = -13740
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19920728
$ws.Range("I4").Value = 22000112
$ws.Range("K4").Value = 66000336
$ws.Range("M4").Value = -66000224
$ws.Range("H94").Value = 10716
$ws.Range("I94").Value = 2972.5715
$ws.Range("K94").Value = 8917.7145
$ws.Range("M94").Value = -8241.7145
$ws.Range("H99").Value = 2133
$ws.Range("I99").Value = 2133
$ws.Range("K99").Value = 6399
$ws.Range("M99").Value = -4153
$ws.Range("H125").Value = 1633.3334
$ws.Range("I125").Value = 1633.3334
$ws.Range("K125").Value = 4900.0002
$ws.Range("M125").Value = 19.9997999999996
$ws.Range("H134").Value = 12258.091
$ws.Range("J134").Value = 18372.715
$ws.Range("L134").Value = 55118.145
$ws.Range("N134").Value = -65258.145
$ws.Range("H140").Value = 6226.8
$ws.Range("I140").Value = 2467.0715
$ws.Range("K140").Value = 7401.2145
$ws.Range("M140").Value = -2221.2145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 397.2857
$ws.Range("I97").Value = 405.16666
$ws.Range("K97").Value = 405.16666
$ws.Range("M97").Value = 90.83334000000002
$ws.Range("H102").Value = 1237.95
$ws.Range("I102").Value = 878.75
$ws.Range("J102").Value = 2674.75
$ws.Range("K102").Value = 878.75
$ws.Range("L102").Value = 2674.75
$ws.Range("M102").Value = 743.25
$ws.Range("N102").Value = -5918.75
$ws.Range("H113").Value = 1096.6
$ws.Range("I113").Value = 1009.375
$ws.Range("K113").Value = 1009.375
$ws.Range("M113").Value = 1160.625
$ws.Range("H122").Value = 3637.7
$ws.Range("I122").Value = 2374.2
$ws.Range("K122").Value = 7122.599999999999
$ws.Range("M122").Value = -4672.599999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1255.5
$ws.Range("I22").Value = 858.8
$ws.Range("K22").Value = 858.8
$ws.Range("M22").Value = -563.8
$ws.Range("H27").Value = 1255.5
$ws.Range("I27").Value = 858.8
$ws.Range("K27").Value = 858.8
$ws.Range("M27").Value = -751.8
$ws.Range("H68").Value = 2483.8333
$ws.Range("J68").Value = 2701.125
$ws.Range("L68").Value = 2701.125
$ws.Range("N68").Value = -4199.125
$ws.Range("H71").Value = 2483.8333
$ws.Range("J71").Value = 2701.125
$ws.Range("L71").Value = 13505.625
$ws.Range("N71").Value = -20993.625
$ws.Range("H93").Value = 2762.7144
$ws.Range("I93").Value = 2386.5557
$ws.Range("J93").Value = 3439.8
$ws.Range("K93").Value = 2386.5557
$ws.Range("L93").Value = 3439.8
$ws.Range("M93").Value = -1138.5557
$ws.Range("N93").Value = -5935.8
$ws.Range("H100").Value = 4322.1113
$ws.Range("I100").Value = 2966.6667
$ws.Range("J100").Value = 4999.8335
$ws.Range("K100").Value = 2966.6667
$ws.Range("L100").Value = 4999.8335
$ws.Range("M100").Value = -2425.6667
$ws.Range("N100").Value = -6081.8335
$ws.Range("H122").Value = 7243.793
$ws.Range("I122").Value = 8552.066000000001
$ws.Range("K122").Value = 25656.198
$ws.Range("M122").Value = -23206.198
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1704.0714
$ws.Range("I132").Value = 1635.6
$ws.Range("K132").Value = 4906.799999999999
$ws.Range("M132").Value = -2376.799999999999
$ws.Range("H136").Value = 2711
$ws.Range("I136").Value = 1584.5
$ws.Range("K136").Value = 4753.5
$ws.Range("M136").Value = -2203.5
